$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "30÷3=10, 0" "53÷8=6, 5"
Replace-Text "94÷4=23, 2" "32÷2=16, 0"
Replace-Text "13÷4=3, 1" "67÷3=22, 1"
Replace-Text "95÷4=23, 3" "25÷7=3, 4"
Replace-Text "79÷7=11, 2" "55÷3=18, 1"
Replace-Text "70÷8=8, 6" "48÷6=8, 0"
Replace-Text "14÷2=7, 0" "62÷9=6, 8"
Replace-Text "17÷3=5, 2" "25÷7=3, 4"
Replace-Text "47÷7=6, 5" "71÷4=17, 3"
Replace-Text "26÷6=4, 2" "79÷6=13, 1"
Replace-Text "56÷8=7, 0" "50÷5=10, 0"
Replace-Text "60÷6=10, 0" "56÷3=18, 2"
Replace-Text "78÷7=11, 1" "75÷2=37, 1"
Replace-Text "83÷4=20, 3" "15÷9=1, 6"
Replace-Text "75÷7=10, 5" "36÷6=6, 0"
Replace-Text "88÷5=17, 3" "96÷2=48, 0"
Replace-Text "63÷6=10, 3" "83÷6=13, 5"
Replace-Text "62÷4=15, 2" "54÷7=7, 5"
Replace-Text "18÷9=2, 0" "78÷2=39, 0"
Replace-Text "27÷5=5, 2" "49÷8=6, 1"
Replace-Text "84÷8=10, 4" "35÷9=3, 8"
Replace-Text "13÷5=2, 3" "78÷9=8, 6"
Replace-Text "25÷2=12, 1" "71÷2=35, 1"
Replace-Text "56÷6=9, 2" "81÷6=13, 3"
Replace-Text "30÷5=6, 0" "77÷4=19, 1"
